$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force the price/volume columns to remain text (matching the source
# workbook's inline-string cells) instead of being auto-converted to
# numbers by Excel's type inference when we assign numeric-looking values.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.906.47"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.361.59"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("E5").Value = "  +5.81%  "
$ws.Range("D6").Value = "241.81"
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("D7").Value = "77.09"
$ws.Range("E7").Value = "  +6.32%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.639"
$ws.Range("E9").Value = "  +28.42%  "
$ws.Range("E10").Value = "  +4.58%  "
$ws.Range("D11").Value = "57.43"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "33.51"
$ws.Range("E12").Value = "  +22.34%  "
$ws.Range("D13").Value = "7.60"
$ws.Range("E13").Value = "  +20.24%  "
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").Value = "2.712.81"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").Value = "0.924"
$ws.Range("E17").Value = "  +6.12%  "
$ws.Range("D18").Value = "2.363.86"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "43.874.91"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  +3.38%  "
$ws.Range("E21").Value = "  +4.65%  "
$ws.Range("D22").Value = "77.67"
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("D23").Value = "257.11"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("D26").Value = "11.14"
$ws.Range("E26").Value = "  +10.00%  "
$ws.Range("D27").Value = "3.61"
$ws.Range("E27").Value = "  -6.90%  "
$ws.Range("E28").Value = "  +15.37%  "
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").Value = "23.15"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").Value = "174.56"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").Value = "0.136"
$ws.Range("E33").Value = "  +6.17%  "
$ws.Range("D34").Value = "5.34"
$ws.Range("E34").Value = "  +6.02%  "
$ws.Range("D35").Value = "0.0757"
$ws.Range("E35").Value = "  +8.81%  "
$ws.Range("D36").Value = "5.40"
$ws.Range("E36").Value = "  +6.85%  "
$ws.Range("D37").Value = "3.82"
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").Value = "0.0278"
$ws.Range("E40").Value = "  +8.36%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "0.204"
$ws.Range("E42").Value = "  +18.96%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "9.01"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("B44").Value = "BinanceUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  +6.27%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "1.27"
$ws.Range("E46").Value = "  +5.27%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").Value = "  +12.97%  "
$ws.Range("E48").Value = "  +2.51%  "
$ws.Range("D49").Value = "102.06"
$ws.Range("E49").Value = "  +2.32%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "55.39"
$ws.Range("E51").Value = "  +8.92%  "

# Drop the temporary text-number-format so the cells end up with the same
# (default/no explicit) style they started with.
$ws.Range("B2:E51").Style = "Normal"
